$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '46.587.49'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  +5.73%  '

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '2.299.28'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  +3.40%  '

$ws.Range('E4').Value = '  -0.12%  '

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '304.74'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +2.17%  '

$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '101.53'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +12.38%  '

$ws.Range('E7').Value = '  +2.05%  '

$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -0.03%  '

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.524'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +6.23%  '

$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '36.60'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  +10.14%  '

$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.0791'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +1.68%  '

$ws.Range('E12').Value = '  +6.95%  '

$ws.Range('E13').Value = '  +0.13%  '

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '2.647.78'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +3.10%  '

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '2.296.14'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +3.13%  '

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '13.84'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +3.46%  '

$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '0.815'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +4.97%  '

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '46.570.04'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +5.92%  '

$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '13.05'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +7.54%  '

$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0944'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +4.14%  '

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '6.02'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +0.68%  '

$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '66.19'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +3.35%  '

$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '249.25'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +6.04%  '

$ws.Range('E24').Value = '  +3.08%  '

$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +0.18%  '

$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '1.93'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +4.96%  '

$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '42.47'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +8.86%  '

$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '2.27'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  +0.18%  '

$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '9.90'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +5.64%  '

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '20.04'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +4.36%  '

$ws.Range('E31').Value = '  +14.09%  '

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '5.65'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +3.09%  '

$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '147.24'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -2.97%  '

$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.0795'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +4.33%  '

$ws.Range('E35').Value = '  +15.09%  '

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.114'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +10.09%  '

$ws.Range('E37').Value = '  +1.07%  '

$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '16.09'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +19.23%  '

$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '1.77'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +5.77%  '

$ws.Range('E40').Value = '  +11.33%  '

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '3.37'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +6.76%  '

$ws.Range('E42').Value = '  +1.26%  '

$ws.Range('E43').Value = '  -0.23%  '

$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '1.99'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +10.93%  '

$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '1.811.84'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +0.88%  '

$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '88.63'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +21.23%  '

$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.196'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +5.49%  '

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '73.32'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +8.11%  '

$ws.Range('E49').Value = '  +6.20%  '

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '95.98'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +1.60%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '2.525.11'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +3.22%  '
